# Rename the field/value mapping sheets for excel/csv consistency,
# and move the active/selected tab from "field_mapping" (now "fields")
# to "value_mapping" (now "values").

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("field_mapping").Name = "fields"
$wb.Worksheets.Item("value_mapping").Name = "values"

$wb.Worksheets.Item("values").Activate()
